$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the three new columns (at their FINAL target letters, applied
#        left-to-right so each successive target already accounts for the
#        columns inserted before it). ---
$ws.Columns("U:U").Insert()   # new "Total Stall Count" column
$ws.Columns("Y:Y").Insert()   # new "Lift Per Tower" column
$ws.Columns("AA:AA").Insert() # new "Flier Frequency" column

# --- 2. Fix up the two edited shared strings (QWE -> ABC, ZXC -> ADF). ---
$ws.Range("E2").Value = "ABC"
$ws.Range("E3").Value = "ADF"

# --- 3. Header row (row 1): set text for the three brand-new header cells
#        and rename the old "Door to Door Flier Allowed(Y/N)" column (which
#        the inserts above shifted to Z1) to "Flier Allowed(Y/N)". Everything
#        to the right (AB1..AF1, the old U1..AC1 headers) just slides along
#        automatically with the column inserts and needs no further edits. ---
$ws.Range("U1").Value = "Total Stall Count"
$ws.Range("Y1").Value = "Lift Per Tower"
$ws.Range("Z1").Value = "Flier Allowed(Y/N)"
$ws.Range("AA1").Value = "Flier Frequency"

# --- 4. Data rows: fill the newly inserted cells with their values. ---
# Row 2
$ws.Range("U2").Value = 1
$ws.Range("Y2").Value = 2
$ws.Range("AA2").Value = 2
# Row 3
$ws.Range("U3").Value = 2
$ws.Range("Y3").Value = 3
$ws.Range("AA3").Value = 3

# --- 5. Sheet view changes: zoom level and active selection. ---
$excel.ActiveWindow.Zoom = 75
$ws.Range("E5").Select() | Out-Null

# --- 6. Row 1 height and column width tweaks (U/Y keep the ~11.52 width of
#        their neighbours, Z - the renamed flier column - narrows slightly,
#        and AA - the new flier-frequency column - takes on the ~12.94 width
#        that the "Stall Price" column used to have). Inputs are chosen so
#        the engine's internal pixel-quantized ColumnWidth lands as close as
#        possible to the target values. ---
$ws.Rows("1:1").RowHeight = 55.75
$ws.Columns("U:U").ColumnWidth = 10.7
$ws.Columns("Y:Y").ColumnWidth = 10.7
$ws.Columns("Z:Z").ColumnWidth = 9.35
$ws.Columns("AA:AA").ColumnWidth = 12.1
